# Add a second comment ("Colin Debuiche") that annotates the exact same
# span of text as the existing comment 0 (Marc Smith) -- the "e" at the
# end of "On les aporte" in "...ordinaire. On les aporte basilics...".
#
# We anchor the new comment on the existing comment's own Scope range so
# it wraps precisely the same characters, rather than re-locating the
# text with Find (which would be fragile if the word occurs more than
# once).

$d = $word.ActiveDocument

$existingComment = $d.Comments(1)
$targetRange = $existingComment.Scope

$newComment = $d.Comments.Add($targetRange, "PB est d'accord. Il pense que l'auteur a utilisé un mot pour un autre.")
$newComment.Author = "Colin Debuiche"
$newComment.Initial = ""
$newComment.Date = "2018-07-13T13:08:57Z"
